# Gate_Planning.xlsx - "Add files via upload" edit
#
# Summary of the change (reconstructed from the OOXML diff):
#  - Gates sheet: column F held a helper list of flight-time windows per
#    gate; all of that data is removed (F2:F29 cleared). Two of the cells
#    (F3, F25) keep their cell-level style (a blank, still-formatted cell)
#    while the rest disappear entirely once cleared.
#  - Flight Schedule sheet: column L held a parallel "gate code" helper
#    column; the whole column is deleted.
#  - Because those helper strings are no longer referenced anywhere, the
#    used range / selections also move; the previously active sheet
#    (Gates) is no longer the active tab - Flight Schedule becomes active.

$wb = $excel.ActiveWorkbook

$wsFlight = $wb.Worksheets.Item("Flight Schedule")
$wsGates  = $wb.Worksheets.Item("Gates")

# --- Gates sheet: drop the helper "time window" column F ---------------
$wsGates.Range("F2:F29").ClearContents()

# --- Flight Schedule sheet: drop the helper gate-code column L ---------
$wsFlight.Columns("L").Delete()

# --- Restore / update selections ----------------------------------------
# Select Gates' old range first (so it is no longer the final/active
# sheet), then finish on Flight Schedule so it becomes the active tab,
# matching the saved view state in the target workbook.
$wsGates.Range("F2:F29").Select()

$wsFlight.Select()
$wsFlight.Range("M13").Select()
